$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns J..Y (10..25), header row 1, plus identical values in rows 2..23
$headers = @(
    "I4 Low Pixel Limit",
    "I4 High Pixel Limit",
    "I2 Low Pixel Limit",
    "I2 High Pixel Limit",
    "I1 Low Pixel Limit",
    "I1 High Pixel Limit",
    "I3 Low Pixel Limit",
    "I3 High Pixel Limit",
    "I4 Contrast",
    "I4 Bias",
    "I2 Contrast",
    "I2 Bias",
    "I1 Contrast",
    "I1 Bias",
    "I3 Contrast",
    "I3 Bias"
)

$values = @(
    3.30368,
    17.5204,
    0.204653,
    15.3575,
    -0.470074,
    15.7846,
    0.431342,
    15.4128,
    2.24647,
    0.508108,
    5.03209,
    0.292973,
    5.26316,
    0.377838,
    3.76866,
    0.425373
)

$firstCol = 10  # column J
$lastDataRow = 23

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $firstCol + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    for ($r = 2; $r -le $lastDataRow; $r++) {
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
